$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oligos")

# Insert a new column before column A ("Name" etc. shift right by one column)
# and use it to hold a short "Tag" for each oligo (derived from the full name).
$ws.Columns("A").Insert()

$ws.Range("A1").Value = "Tag"
$ws.Range("A2").Value = "o2"
$ws.Range("A3").Value = "o3"
$ws.Range("A4").Value = "o4"

# The new Tag column (and the Name column next to it) are narrow, short
# identifiers, so size them down from the old wide "Name" column width.
$ws.Range("A1:B1").ColumnWidth = 8.95

# Update the print area defined name so it still points at the same logical
# cells now that everything shifted one column to the right.
$ws.PageSetup.PrintArea = '$I$94:$J$98'

# Leave selection on the row below the data, first column.
$ws.Range("A5").Select() | Out-Null
